$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 39996.668
$ws.Range("J109").Value = 39996.668
$ws.Range("L109").Value = 39996.668
$ws.Range("N109").Value = -42770.668

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H129").Value = 1735.0728
$ws.Range("J129").Value = 1563.04
$ws.Range("L129").Value = 4689.12
$ws.Range("N129").Value = -14689.12

$ws.Range("H133").Value = 57266.75
$ws.Range("J133").Value = 57266.75
$ws.Range("L133").Value = 57266.75
$ws.Range("N133").Value = -67386.75

$ws.Range("H137").Value = 4168.587
$ws.Range("I137").Value = 1454.6451
$ws.Range("J137").Value = 9777.4
$ws.Range("K137").Value = 4363.9353
$ws.Range("L137").Value = 29332.2
$ws.Range("M137").Value = -1813.9353
$ws.Range("N137").Value = -34432.2

$ws.Range("H138").Value = 2482.4526
$ws.Range("I138").Value = 1411.5714
$ws.Range("J138").Value = 3331.0754
$ws.Range("K138").Value = 4234.7142
$ws.Range("L138").Value = 9993.226200000001
$ws.Range("M138").Value = 905.2857999999997
$ws.Range("N138").Value = -20273.2262

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1659.6794
$ws.Range("I74").Value = 1539.9839
$ws.Range("J74").Value = 2123.5
$ws.Range("K74").Value = 1539.9839
$ws.Range("L74").Value = 2123.5
$ws.Range("M74").Value = -665.9838999999999
$ws.Range("N74").Value = -3871.5

$ws.Range("H75").Value = 24000
$ws.Range("J75").Value = 24000
$ws.Range("L75").Value = 24000
$ws.Range("N75").Value = -25748

$ws.Range("H77").Value = 1659.6794
$ws.Range("I77").Value = 1539.9839
$ws.Range("J77").Value = 2123.5
$ws.Range("K77").Value = 7699.9195
$ws.Range("L77").Value = 10617.5
$ws.Range("M77").Value = -3331.9195
$ws.Range("N77").Value = -19353.5

$ws.Range("H78").Value = 24000
$ws.Range("J78").Value = 24000
$ws.Range("L78").Value = 72000
$ws.Range("N78").Value = -80736

$ws.Range("H117").Value = 48412.668
$ws.Range("J117").Value = 48412.668
$ws.Range("L117").Value = 48412.668
$ws.Range("N117").Value = -57590.668

$ws.Range("H118").Value = 49202
$ws.Range("J118").Value = 49202
$ws.Range("L118").Value = 49202
$ws.Range("N118").Value = -52516

$ws.Range("H132").Value = 8066207.5
$ws.Range("I132").Value = 12821719
$ws.Range("J132").Value = 2513.913
$ws.Range("K132").Value = 38465157
$ws.Range("L132").Value = 7541.739
$ws.Range("M132").Value = -38462627
$ws.Range("N132").Value = -12601.739

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 44798.2
$ws.Range("J133").Value = 44798.2
$ws.Range("L133").Value = 44798.2
$ws.Range("N133").Value = -54918.2

$ws.Range("H134").Value = 145949.97
$ws.Range("I134").Value = 1356.5385
$ws.Range("J134").Value = 198164.27
$ws.Range("K134").Value = 4069.6155
$ws.Range("L134").Value = 594492.8099999999
$ws.Range("M134").Value = -1534.6155
$ws.Range("N134").Value = -599562.8099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 659.1818
$ws.Range("I7").Value = 668.875
$ws.Range("J7").Value = 633.3333
$ws.Range("K7").Value = 668.875
$ws.Range("L7").Value = 633.3333
$ws.Range("M7").Value = -555.875
$ws.Range("N7").Value = -859.3333

$ws.Range("H19").Value = 225
$ws.Range("I19").Value = 225
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 225
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -55
$ws.Range("N19").ClearContents()

$ws.Range("H22").Value = 2006.2307
$ws.Range("I22").Value = 422.33334
$ws.Range("J22").Value = 5570
$ws.Range("K22").Value = 422.33334
$ws.Range("L22").Value = 5570
$ws.Range("M22").Value = -72.33334000000002
$ws.Range("N22").Value = -6270

$ws.Range("H24").Value = 225
$ws.Range("I24").Value = 225
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 225
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -55
$ws.Range("N24").ClearContents()

$ws.Range("H31").Value = 2975.04
$ws.Range("I31").Value = 1320.037
$ws.Range("J31").Value = 3587.1643
$ws.Range("K31").Value = 1320.037
$ws.Range("L31").Value = 3587.1643
$ws.Range("M31").Value = -1025.037
$ws.Range("N31").Value = -4177.1643

$ws.Range("H34").Value = 2975.04
$ws.Range("I34").Value = 1320.037
$ws.Range("J34").Value = 3587.1643
$ws.Range("K34").Value = 1320.037
$ws.Range("L34").Value = 3587.1643
$ws.Range("M34").Value = -1118.037
$ws.Range("N34").Value = -3991.1643

$ws.Range("H122").Value = 134377.56
$ws.Range("I122").Value = 240700.6
$ws.Range("K122").Value = 722101.8
$ws.Range("M122").Value = -719651.8

$ws.Range("H132").Value = 50099.9
$ws.Range("I132").Value = 1660.8695
$ws.Range("J132").Value = 235782.83
$ws.Range("K132").Value = 4982.6085
$ws.Range("L132").Value = 707348.49
$ws.Range("M132").Value = -2452.6085
$ws.Range("N132").Value = -712408.49

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3170.3809
$ws.Range("I5").Value = 4329.077
$ws.Range("J5").Value = 1287.5
$ws.Range("K5").Value = 12987.231
$ws.Range("L5").Value = 3862.5
$ws.Range("M5").Value = -12875.231
$ws.Range("N5").Value = -4086.5

$ws.Range("H113").Value = 3655.7878
$ws.Range("I113").Value = 5867.9473
$ws.Range("J113").Value = 653.5714
$ws.Range("K113").Value = 17603.8419
$ws.Range("L113").Value = 1960.7142
$ws.Range("M113").Value = -15433.8419
$ws.Range("N113").Value = -6300.7142

$ws.Range("H132").Value = 2128.7856
$ws.Range("I132").Value = 950.5
$ws.Range("J132").Value = 3699.8333
$ws.Range("K132").Value = 8554.5
$ws.Range("L132").Value = 33298.4997
$ws.Range("M132").Value = -6024.5
$ws.Range("N132").Value = -38358.4997

$ws.Range("H135").Value = 3170.3809
$ws.Range("I135").Value = 4329.077
$ws.Range("J135").Value = 1287.5
$ws.Range("K135").Value = 38961.693
$ws.Range("L135").Value = 11587.5
$ws.Range("M135").Value = -36426.693
$ws.Range("N135").Value = -16657.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1230.0714
$ws.Range("I122").Value = 1265.5454
$ws.Range("K122").Value = 3796.6362
$ws.Range("M122").Value = -1346.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2247.8445
$ws.Range("I132").Value = 1620.9642
$ws.Range("J132").Value = 3280.353
$ws.Range("K132").Value = 4862.892599999999
$ws.Range("L132").Value = 9841.059000000001
$ws.Range("M132").Value = -2332.892599999999
$ws.Range("N132").Value = -14901.059

$ws.Range("H136").Value = 1340.3846
$ws.Range("I136").Value = 1054.6364
$ws.Range("K136").Value = 3163.9092
$ws.Range("M136").Value = -613.9092000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 73060.5
$ws.Range("J46").Value = 73060.5
$ws.Range("L46").Value = 73060.5
$ws.Range("N46").Value = -73522.5

$ws.Range("H122").Value = 9524809
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 1260.5541
$ws.Range("I132").Value = 1161.7234
$ws.Range("J132").Value = 1432.5927
$ws.Range("K132").Value = 3485.1702
$ws.Range("L132").Value = 4297.7781
$ws.Range("M132").Value = -955.1702000000005
$ws.Range("N132").Value = -9357.7781

$ws.Range("H134").Value = 73060.5
$ws.Range("J134").Value = 73060.5
$ws.Range("L134").Value = 219181.5
$ws.Range("N134").Value = -224251.5

Write-Host "Applied all leve profit updates across 8 sheets"
